# This edit permutes the data held in columns D and K:T across rows 2-13
# (everything else - A,B,C,E,F,G,H,I,J - is identical on every row so it
# is left untouched). The mapping of "row now holds the data that used to
# live on row X" is:
#   2<-3, 3<-4, 4<-5, 5<-2, 6<-13, 7<-11, 8<-10, 9<-8, 10<-9, 11<-6, 12<-7, 13<-12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 3
    3  = 4
    4  = 5
    5  = 2
    6  = 13
    7  = 11
    8  = 10
    9  = 8
    10 = 9
    11 = 6
    12 = 7
    13 = 12
}

# Columns that move together as a row's "data" (date, variety, quality,
# volume, min/max/avg price, unit, origin, $/kg, kg/unit).
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot every source row's values BEFORE any writes happen, since the
# permutation reads from rows that will themselves be overwritten.
$snapshot = @{}
foreach ($r in 2..13) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write each destination row with the snapshot taken from its source row.
foreach ($destRow in 2..13) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
